$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update financial figures to restated values
$ws.Range("D2").Value = 135
$ws.Range("E2").Value = -17
$ws.Range("F2").Value = -14
$ws.Range("G2").Value = -34
$ws.Range("H2").Value = -35
$ws.Range("I2").Value = -35
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 282
$ws.Range("L2").Value = 97
$ws.Range("M2").Value = 185
$ws.Range("N2").Value = 181
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 121
$ws.Range("Q2").Value = -22
$ws.Range("R2").Value = -36
$ws.Range("S2").Value = 70
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = -22
$ws.Range("V2").Value = 72
$ws.Range("W2").Value = -12.9
$ws.Range("X2").Value = -25.88
$ws.Range("Y2").Value = -18.52
$ws.Range("Z2").Value = -12.73
$ws.Range("AA2").Value = 52.48
$ws.Range("AB2").Value = 51.34
$ws.Range("AC2").Value = -170
$ws.Range("AD2").Value = -3.71
$ws.Range("AE2").Value = 747
$ws.Range("AF2").Value = 0.84
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 24244800

# Row 3: update financial figures to restated values
$ws.Range("D3").Value = 213
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 11
$ws.Range("H3").Value = 11
$ws.Range("I3").Value = 11
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 402
$ws.Range("L3").Value = 109
$ws.Range("M3").Value = 293
$ws.Range("N3").Value = 285
$ws.Range("O3").Value = 8
$ws.Range("P3").Value = 185
$ws.Range("Q3").Value = -22
$ws.Range("R3").Value = -55
$ws.Range("S3").Value = 81
$ws.Range("T3").Value = 5
$ws.Range("U3").Value = -26
$ws.Range("V3").Value = 63
$ws.Range("W3").Value = 4.75
$ws.Range("X3").Value = 5.17
$ws.Range("Y3").Value = 4.74
$ws.Range("Z3").Value = 3.21
$ws.Range("AA3").Value = 37.19
$ws.Range("AB3").Value = 55.73
$ws.Range("AC3").Value = 36
$ws.Range("AD3").Value = 27.27
$ws.Range("AE3").Value = 817
$ws.Range("AF3").Value = 1.19
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 36960027

# Row 4: update financial figures to restated values
$ws.Range("D4").Value = 185
$ws.Range("E4").Value = -9
$ws.Range("F4").Value = -9
$ws.Range("G4").Value = -15
$ws.Range("H4").Value = -15
$ws.Range("I4").Value = -15
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 398
$ws.Range("L4").Value = 55
$ws.Range("M4").Value = 344
$ws.Range("N4").Value = 332
$ws.Range("O4").Value = 12
$ws.Range("P4").Value = 224
$ws.Range("Q4").Value = -8
$ws.Range("R4").Value = -12
$ws.Range("S4").Value = 18
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = -9
$ws.Range("V4").Value = 8
$ws.Range("W4").Value = -4.63
$ws.Range("X4").Value = -8.27
$ws.Range("Y4").Value = -4.86
$ws.Range("Z4").Value = -3.82
$ws.Range("AA4").Value = 15.86
$ws.Range("AB4").Value = 50.11
$ws.Range("AC4").Value = -38
$ws.Range("AD4").Value = -22.07
$ws.Range("AE4").Value = 778
$ws.Range("AF4").Value = 1.09
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 44787538

# Row 5: update financial figures to restated values
$ws.Range("D5").Value = 217
$ws.Range("E5").Value = -29
$ws.Range("F5").Value = -29
$ws.Range("G5").Value = -112
$ws.Range("H5").Value = -112
$ws.Range("I5").Value = -108
$ws.Range("J5").Value = -3
$ws.Range("K5").Value = 383
$ws.Range("L5").Value = 162
$ws.Range("M5").Value = 220
$ws.Range("N5").Value = 223
$ws.Range("O5").Value = -2
$ws.Range("P5").Value = 224
$ws.Range("Q5").Value = 9
$ws.Range("R5").Value = -76
$ws.Range("S5").Value = 103
$ws.Range("T5").Value = 43
$ws.Range("U5").Value = -34
$ws.Range("V5").Value = 103
$ws.Range("W5").Value = -13.35
$ws.Range("X5").Value = -51.5
$ws.Range("Y5").Value = -39.07
$ws.Range("Z5").Value = -28.65
$ws.Range("AA5").Value = 73.5
$ws.Range("AB5").Value = 0.97
$ws.Range("AC5").Value = -242
$ws.Range("AD5").Value = -2.96
$ws.Range("AE5").Value = 522
$ws.Range("AF5").Value = 1.38
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 44787538

# Row 6: update financial figures to restated values
$ws.Range("D6").Value = 323
$ws.Range("E6").Value = -20
$ws.Range("F6").Value = -20
$ws.Range("G6").Value = -31
$ws.Range("H6").Value = -30
$ws.Range("I6").Value = -25
$ws.Range("K6").Value = 430
$ws.Range("L6").Value = 187
$ws.Range("M6").Value = 242
$ws.Range("N6").Value = 231
$ws.Range("P6").Value = 271
$ws.Range("Q6").Value = -5
$ws.Range("R6").Value = -14
$ws.Range("S6").Value = 30
$ws.Range("T6").Value = 9
$ws.Range("U6").Value = -14
$ws.Range("V6").Value = 75
$ws.Range("W6").Value = -6.1
$ws.Range("X6").Value = -9.380000000000001
$ws.Range("Y6").Value = -10.97
$ws.Range("Z6").Value = -7.46
$ws.Range("AA6").Value = 77.3
$ws.Range("AB6").Value = -9.98
$ws.Range("AC6").Value = -52
$ws.Range("AD6").Value = -9.68
$ws.Range("AE6").Value = 443
$ws.Range("AF6").Value = 1.14
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 54212316

# Rows 7-9: remove financial data cells, keeping only index/category/year columns
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()

